$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"
$ws.Range("C3").Value = "Wheat."
$ws.Range("C4").Value = "Wheat."
$ws.Range("D4").Value = "Zorba's Tzatziki, BCfresh Tomatoes"

[void]$ws.Range("E11").Select()
